# Generate Report for Handback
#
# This script updates the localization-status workbook to reflect a
# completed handback: the zh-cn / de-de "Latest Target File" and
# "Latest Handback File" columns get populated, the de-de handback
# timestamp is refreshed, the Overview/Status text changes from
# "Ready for handoff" to "Handed back: in sync with en-US", and the
# widened Status-ish columns get resized.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$mdDisplay  = "cf56f4a4-a726-4997-95b5-bbb4ee2810ce.md"
$mdTarget   = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d08ada96d89118104503fe0ad5a253a9d590390d/e2e/cf56f4a4-a726-4997-95b5-bbb4ee2810ce.md"
$zhXlf      = "cf56f4a4-a726-4997-95b5-bbb4ee2810ce.0cf7c83772e713ab7875ea74fda76374b4389efa.zh-cn.xlf"
$deXlf      = "cf56f4a4-a726-4997-95b5-bbb4ee2810ce.0cf7c83772e713ab7875ea74fda76374b4389efa.de-de.xlf"
$zhHandbackDate = "2016-08-18 15:09:24"
$deHandbackDate = "2016-08-18 15:09:33"

$statusText = "Handed back: in sync with en-US"

# ---- zh-cn sheet: rows 2 and 3 ----
foreach ($row in 2, 3) {
    $zhcn.Range("C$row").Value2 = $statusText
    $zhcn.Range("I$row").Value2 = $mdDisplay
    $zhcn.Hyperlinks.Add($zhcn.Range("I$row"), $mdTarget, "", "", $mdDisplay) | Out-Null
    $zhcn.Range("J$row").Value2 = $zhXlf
    $zhcn.Range("K$row").Value2 = $zhHandbackDate
}

# ---- de-de sheet: rows 2 and 3 ----
foreach ($row in 2, 3) {
    $dede.Range("C$row").Value2 = $statusText
    $dede.Range("I$row").Value2 = $mdDisplay
    $dede.Hyperlinks.Add($dede.Range("I$row"), $mdTarget, "", "", $mdDisplay) | Out-Null
    $dede.Range("J$row").Value2 = $deXlf
    $dede.Range("K$row").Value2 = $deHandbackDate
}

# ---- Overview sheet: Status text ----
$overview = $wb.Worksheets.Item("Overview")
foreach ($row in 2, 3) {
    $overview.Range("E$row").Value2 = $statusText
    $overview.Range("F$row").Value2 = $statusText
}

# ---- Column widths ----
$overview.Columns.Item(5).ColumnWidth = 29.9777047293527
$overview.Columns.Item(6).ColumnWidth = 29.9777047293527

foreach ($ws in $zhcn, $dede) {
    $ws.Columns.Item(3).ColumnWidth = 29.9777047293527
    $ws.Columns.Item(9).ColumnWidth = 40
    $ws.Columns.Item(10).ColumnWidth = 40
}
